$d = $word.ActiveDocument

# 1. Title / headline text (appears as Heading1 at top and as bold text near the end)
$d.Content.Find.Execute(
    "Play Football Mania Deluxe Free - Exciting Football-Themed Online Slot",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Football Mania Deluxe for Free", 2)

# 2. "What we like" bullet list updates
$d.Content.Find.Execute(
    "Exciting football-themed environment",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Thrilling online slot game based on football", 2)

$d.Content.Find.Execute(
    "Innovative gaming features",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exciting special features that enhance gameplay", 2)

$d.Content.Find.Execute(
    "Special football-related bonus features",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sporty environment with impressive graphics", 2)

$d.Content.Find.Execute(
    "Customizable game experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Innovative gaming features for customization", 2)

# 3. "What we don't like" bullet list updates
$d.Content.Find.Execute(
    "Limited number of reels",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited number of reels and paylines", 2)

# 4. Insert a new bullet paragraph after "Limited number of reels and paylines",
#    matching the same ListBullet style/formatting.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Limited number of reels and paylines") {
        $target = $p
        break
    }
}
$newPara = $target.Range.InsertParagraphAfter()
$newRange = $target.Next().Range
$newRange.Text = "May not appeal to players who are not fans of football"

# 5. Closing bold headline (same text as title, already handled by step 1's
#    whole-document Find/Replace) and the closing italic summary text.
$d.Content.Find.Execute(
    "Join the team and try to become the best player. Play Football Mania Deluxe free and enjoy exciting football-related bonus features in an innovative customizable game experience.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Football Mania Deluxe and play this exciting online slot game for free.", 2)
